$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 11 (years 2000年-2009年), shifting the remaining
# rows (2010年-2015年) up so they become rows 2-7.
$ws.Range("A2:D11").EntireRow.Delete()
